$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold a literal text value (no numeric/string
# auto-coercion by Excel, and no lingering custom style/number-format left
# behind on the cell once we're done).
function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "27.794.03"
$ws.Range("E2").Value = "  +1.63%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.879.35"
$ws.Range("E3").Value = "  +1.26%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "331.69"
$ws.Range("E5").Value = "  +2.60%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.09%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4731"
$ws.Range("E7").Value = "  +4.66%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3955"
$ws.Range("E8").Value = "  +2.32%  "

# Row 9 - OKB
Set-TextValue "D9" "48.00"
$ws.Range("E9").Value = "  -0.43%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.08085"
$ws.Range("E10").Value = "  +2.33%  "

# Row 11 - Polygon
Set-TextValue "D11" "1.031"
$ws.Range("E11").Value = "  +1.54%  "

# Row 12 - Solana
$ws.Range("E12").Value = "  +3.95%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.899.73"
$ws.Range("E13").Value = "  +2.15%  "

# Row 14 - Polkadot
Set-TextValue "D14" "5.970"
$ws.Range("E14").Value = "  +1.10%  "

# Row 15 - Chainlink
Set-TextValue "D15" "7.148"
$ws.Range("E15").Value = "  +0.56%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  +0.13%  "

# Row 17 - was Litecoin, now ShibaInu
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.00001051"
$ws.Range("E17").Value = "  +2.14%  "

# Row 18 - was ShibaInu, now Litecoin
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D18" "87.15"
$ws.Range("E18").Value = "  +1.58%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06668"
$ws.Range("E19").Value = "  +1.82%  "

# Row 20 - Avalanche
Set-TextValue "D20" "17.31"
$ws.Range("E20").Value = "  +1.82%  "

# Row 22 - WrappedBTC
Set-TextValue "D22" "27.805.46"
$ws.Range("E22").Value = "  +1.70%  "

# Row 23 - Uniswap
Set-TextValue "D23" "5.535"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24 - Cosmos
Set-TextValue "D24" "11.01"
$ws.Range("E24").Value = "  +0.51%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +1.00%  "

# Row 26 - WrappedliquidstakedEther2.0
Set-TextValue "D26" "2.106.73"
$ws.Range("E26").Value = "  +1.23%  "

# Row 27 - Monero
Set-TextValue "D27" "159.07"
$ws.Range("E27").Value = "  +3.53%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "20.23"
$ws.Range("E28").Value = "  +1.69%  "

# Row 29 - LidoDAOToken
Set-TextValue "D29" "2.112"
$ws.Range("E29").Value = "  +1.82%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "5.614"
$ws.Range("E30").Value = "  +3.38%  "

# Row 31 - BitcoinCash
Set-TextValue "D31" "122.48"
$ws.Range("E31").Value = "  +1.28%  "

# Row 32 - ImmutableX
Set-TextValue "D32" "0.9913"
$ws.Range("E32").Value = "  +5.95%  "

# Row 33 - Stellar
Set-TextValue "D33" "0.09560"
$ws.Range("E33").Value = "  +3.13%  "

# Row 34 - ARBITRUM
Set-TextValue "D34" "1.455"
$ws.Range("E34").Value = "  -1.59%  "

# Row 35 - HuobiToken
Set-TextValue "D35" "3.590"
$ws.Range("E35").Value = "  -0.28%  "

# Row 36 - Filecoin
Set-TextValue "D36" "5.358"
$ws.Range("E36").Value = "  +1.78%  "

# Row 37 - Hedera
Set-TextValue "D37" "0.06119"
$ws.Range("E37").Value = "  +2.21%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.02261"
$ws.Range("E38").Value = "  +1.56%  "

# Row 39 - TrustWalletToken
Set-TextValue "D39" "1.236"
$ws.Range("E39").Value = "  +0.18%  "

# Row 40 - FraxShare
Set-TextValue "D40" "8.140"
$ws.Range("E40").Value = "  -0.35%  "

# Row 41 - TheSandbox
Set-TextValue "D41" "0.6038"
$ws.Range("E41").Value = "  +2.31%  "

# Row 42 - Algorand
Set-TextValue "D42" "0.1909"
$ws.Range("E42").Value = "  +0.67%  "

# Row 43 - Aptos
Set-TextValue "D43" "10.29"
$ws.Range("E43").Value = "  +1.68%  "

# Row 44 - Decentraland
Set-TextValue "D44" "0.5740"
$ws.Range("E44").Value = "  +2.17%  "

# Row 45 - WEMIXTOKEN
Set-TextValue "D45" "1.258"
$ws.Range("E45").Value = "  -1.28%  "

# Row 46 - EnergySwap
Set-TextValue "D46" "12.21"
$ws.Range("E46").Value = "  +1.20%  "

# Row 47 - NEARProtocol
Set-TextValue "D47" "1.950"
$ws.Range("E47").Value = "  +1.75%  "

# Row 48 - PancakeSwap
Set-TextValue "D48" "3.381"
$ws.Range("E48").Value = "  +0.22%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.06899"
$ws.Range("E49").Value = "  +2.00%  "

# Row 50 - Quant
Set-TextValue "D50" "114.74"
$ws.Range("E50").Value = "  +5.93%  "

# Row 51 - EOS
Set-TextValue "D51" "1.074"
$ws.Range("E51").Value = "  +2.11%  "
